$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data1 = New-Object 'object[,]' 12,10
$data1[0,0] = 0.07829464001795738
$data1[0,1] = 0.07877937800938027
$data1[0,2] = 0.07864014179500259
$data1[0,3] = 0.07845946774641657
$data1[0,4] = 0.07839212744846909
$data1[0,5] = 0.07822152187813221
$data1[0,6] = 0.07854383010459076
$data1[0,7] = 0.07873393418916634
$data1[0,8] = 0.07847923265723404
$data1[0,9] = 0.0789098965158269
$data1[1,0] = 25.85019593345956
$data1[1,1] = 25.85070029556539
$data1[1,2] = 25.85099975479828
$data1[1,3] = 25.84861473840749
$data1[1,4] = 25.84971866607145
$data1[1,5] = 25.85125078759752
$data1[1,6] = 25.85014495677666
$data1[1,7] = 25.84863092719272
$data1[1,8] = 25.85136409568195
$data1[1,9] = 25.8508551500793
$data1[2,0] = 151.6999346847188
$data1[2,1] = 151.7000764647251
$data1[2,2] = 151.7002375642224
$data1[2,3] = 151.7001719018875
$data1[2,4] = 151.7003755474327
$data1[2,5] = 151.7003002664522
$data1[2,6] = 151.6996346577992
$data1[2,7] = 151.700029767051
$data1[2,8] = 151.7001880554535
$data1[2,9] = 151.6998287745109
$data1[3,0] = 0.1860831622481827
$data1[3,1] = 0.1862753780545544
$data1[3,2] = 0.1857153692294063
$data1[3,3] = 0.1858021064205384
$data1[3,4] = 0.1859662790003875
$data1[3,5] = 0.1858419944894474
$data1[3,6] = 0.1862059983888287
$data1[3,7] = 0.185682111945499
$data1[3,8] = 0.185744846220143
$data1[3,9] = 0.1862053133238964
$data1[4,0] = 2.26001763044322
$data1[4,1] = 2.260004186377163
$data1[4,2] = 2.260054718534799
$data1[4,3] = 2.260063664649345
$data1[4,4] = 2.25999867772061
$data1[4,5] = 2.259959387229524
$data1[4,6] = 2.259987328146046
$data1[4,7] = 2.260002012944936
$data1[4,8] = 2.260034208203687
$data1[4,9] = 2.260024586863321
$data1[5,0] = 3.453675847372491
$data1[5,1] = 3.455541528031846
$data1[5,2] = 3.457999755294844
$data1[5,3] = 3.456235084187883
$data1[5,4] = 3.452212508390573
$data1[5,5] = 3.457980330006651
$data1[5,6] = 3.451763697386524
$data1[5,7] = 3.464101212240302
$data1[5,8] = 3.452743040987104
$data1[5,9] = 3.447870103102793
$data1[6,0] = 28.34962624529889
$data1[6,1] = 28.344000833388
$data1[6,2] = 28.34835403664405
$data1[6,3] = 28.35635819023252
$data1[6,4] = 28.3519181595647
$data1[6,5] = 28.34910650015055
$data1[6,6] = 28.34433914755335
$data1[6,7] = 28.35180112008402
$data1[6,8] = 28.35099152884476
$data1[6,9] = 28.35002394003423
$data1[7,0] = 1.019583450301445
$data1[7,1] = 0.9822501316885365
$data1[7,2] = 0.9261442989850726
$data1[7,3] = 1.047552846343779
$data1[7,4] = 0.9520601549159953
$data1[7,5] = 1.006742665056656
$data1[7,6] = 1.044927020788788
$data1[7,7] = 0.9707211411912515
$data1[7,8] = 0.9545144623237745
$data1[7,9] = 1.084428070334901
$data1[8,0] = 1.007534181037856
$data1[8,1] = 0.9546461279385999
$data1[8,2] = 0.9719281590415387
$data1[8,3] = 0.9721441075529493
$data1[8,4] = 1.046107470645524
$data1[8,5] = 1.034610389099874
$data1[8,6] = 0.9371348491067988
$data1[8,7] = 0.9632213719002793
$data1[8,8] = 1.072636439813636
$data1[8,9] = 0.9336179197863443
$data1[9,0] = 2.148271484035093
$data1[9,1] = 2.147482920061125
$data1[9,2] = 2.148346146262153
$data1[9,3] = 2.147891282509959
$data1[9,4] = 2.148129718447739
$data1[9,5] = 2.148097846148703
$data1[9,6] = 2.147756030947023
$data1[9,7] = 2.148159141364875
$data1[9,8] = 2.147436942348215
$data1[9,9] = 2.147558505758274
$data1[10,0] = 0.9139772165975671
$data1[10,1] = 0.9141457440875383
$data1[10,2] = 0.9138312621800739
$data1[10,3] = 0.9139871782117261
$data1[10,4] = 0.914051946331627
$data1[10,5] = 0.914302973228477
$data1[10,6] = 0.9140506423466492
$data1[10,7] = 0.9142077091995608
$data1[10,8] = 0.9138909900399609
$data1[10,9] = 0.9143307422070915
$data1[11,0] = 0.07853833481992255
$data1[11,1] = 0.07861817272536452
$data1[11,2] = 0.07841600278424737
$data1[11,3] = 0.07777917458271301
$data1[11,4] = 0.07871781409027048
$data1[11,5] = 0.07868380093850569
$data1[11,6] = 0.07878864845778275
$data1[11,7] = 0.07886209056618407
$data1[11,8] = 0.07879170469206687
$data1[11,9] = 0.07831237626461238
$ws.Range("B2:K13").Value = $data1

$data2 = New-Object 'object[,]' 18,10
$data2[0,0] = 0.5573914293617266
$data2[0,1] = 0.5574819437438194
$data2[0,2] = 0.5559152672228114
$data2[0,3] = 0.5578127642838832
$data2[0,4] = 0.557194318724387
$data2[0,5] = 0.5560378759442558
$data2[0,6] = 0.5574108306746501
$data2[0,7] = 0.5600151223014637
$data2[0,8] = 0.5593414678296798
$data2[0,9] = 0.5547791093697644
$data2[1,0] = 0.1765566216305918
$data2[1,1] = 0.1790828072834728
$data2[1,2] = 0.1768632443268893
$data2[1,3] = 0.1755092347525782
$data2[1,4] = 0.1813022432812508
$data2[1,5] = 0.1755682698221086
$data2[1,6] = 0.1786419316875567
$data2[1,7] = 0.1787515149482427
$data2[1,8] = 0.1793819466262752
$data2[1,9] = 0.1781278540065145
$data2[2,0] = 0.1287221680001877
$data2[2,1] = 0.1255264507806259
$data2[2,2] = 0.1259292100341978
$data2[2,3] = 0.1251207486136789
$data2[2,4] = 0.1292148721073986
$data2[2,5] = 0.1267313761644601
$data2[2,6] = 0.1243981675619309
$data2[2,7] = 0.1311479043933923
$data2[2,8] = 0.125248699610824
$data2[2,9] = 0.1283314698565677
$data2[3,0] = 5.105778973152722
$data2[3,1] = 5.109442297168108
$data2[3,2] = 5.108964955443522
$data2[3,3] = 5.107837790213064
$data2[3,4] = 5.113076611127469
$data2[3,5] = 5.107342658957783
$data2[3,6] = 5.106851266750354
$data2[3,7] = 5.108736476894448
$data2[3,8] = 5.110864059677319
$data2[3,9] = 5.106624826647602
$data2[4,0] = 0.03498530696886858
$data2[4,1] = 0.04610574957989958
$data2[4,2] = 0.01811912589879062
$data2[4,3] = 0.19490294959918
$data2[4,4] = 0.0119782390300526
$data2[4,5] = -0.1735601647493568
$data2[4,6] = 0.1159149686986259
$data2[4,7] = 0.008905986893584577
$data2[4,8] = -0.02052641133729731
$data2[4,9] = -0.1335436431116351
$data2[5,0] = 0.4717701502357043
$data2[5,1] = 0.4673506793146259
$data2[5,2] = 0.4696720409254531
$data2[5,3] = 0.4677877521550363
$data2[5,4] = 0.4690655900286287
$data2[5,5] = 0.4688643044983773
$data2[5,6] = 0.4676974214842532
$data2[5,7] = 0.4704527467522187
$data2[5,8] = 0.4662722521868444
$data2[5,9] = 0.4676686963059437
$data2[6,0] = 28.34777339226817
$data2[6,1] = 28.35274637830266
$data2[6,2] = 28.35219107947395
$data2[6,3] = 28.35472398239264
$data2[6,4] = 28.34942630968719
$data2[6,5] = 28.35337801944043
$data2[6,6] = 28.34793007244464
$data2[6,7] = 28.35152330430449
$data2[6,8] = 28.34854285970777
$data2[6,9] = 28.35073385585933
$data2[7,0] = 29.29981278272576
$data2[7,1] = 29.30003614354623
$data2[7,2] = 29.29969333681713
$data2[7,3] = 29.29936336139002
$data2[7,4] = 29.29979820362986
$data2[7,5] = 29.29973687814825
$data2[7,6] = 29.29955298175352
$data2[7,7] = 29.30012856157847
$data2[7,8] = 29.30003950881936
$data2[7,9] = 29.30002670310249
$data2[8,0] = 0.08224020047053716
$data2[8,1] = 0.08194636649499315
$data2[8,2] = 0.08186422590092685
$data2[8,3] = 0.08164065438815642
$data2[8,4] = 0.0818511814526624
$data2[8,5] = 0.08179044892640952
$data2[8,6] = 0.08182374031526424
$data2[8,7] = 0.08154988090471246
$data2[8,8] = 0.08155819761380108
$data2[8,9] = 0.08180848758817866
$data2[9,0] = 0.454874150029585
$data2[9,1] = 0.4545835476838459
$data2[9,2] = 0.4546264619438994
$data2[9,3] = 0.4545633222288285
$data2[9,4] = 0.4547887067399076
$data2[9,5] = 0.4549442585330364
$data2[9,6] = 0.4546211175806204
$data2[9,7] = 0.4547563877505353
$data2[9,8] = 0.4544312934530398
$data2[9,9] = 0.4549797119701208
$data2[10,0] = 1.705118640554433
$data2[10,1] = 1.705723267727757
$data2[10,2] = 1.707485448234782
$data2[10,3] = 1.705621248318673
$data2[10,4] = 1.70682642231767
$data2[10,5] = 1.705593589127088
$data2[10,6] = 1.706185120607468
$data2[10,7] = 1.706366849841732
$data2[10,8] = 1.70677684862712
$data2[10,9] = 1.706176330402858
$data2[11,0] = 3.597729040931274
$data2[11,1] = 3.597175614431924
$data2[11,2] = 3.593129072449697
$data2[11,3] = 3.595088140299688
$data2[11,4] = 3.598501945693067
$data2[11,5] = 3.597435947351289
$data2[11,6] = 3.596697499929983
$data2[11,7] = 3.5909230147504
$data2[11,8] = 3.598891063330766
$data2[11,9] = 3.596479512886261
$data2[12,0] = 12.066180680864
$data2[12,1] = 12.06544430953787
$data2[12,2] = 12.06510201034599
$data2[12,3] = 12.06534929841172
$data2[12,4] = 12.06598123951287
$data2[12,5] = 12.06605208944585
$data2[12,6] = 12.06489323244117
$data2[12,7] = 12.06586595339339
$data2[12,8] = 12.06681303230606
$data2[12,9] = 12.06605055274717
$data2[13,0] = 33.86776796434634
$data2[13,1] = 34.1073979171316
$data2[13,2] = 33.47609118997086
$data2[13,3] = 33.58289528355159
$data2[13,4] = 33.21084359972032
$data2[13,5] = 34.2432549357861
$data2[13,6] = 34.08796274080287
$data2[13,7] = 33.45285365052489
$data2[13,8] = 33.24261212657254
$data2[13,9] = 33.36581112764438
$data2[14,0] = 73826.62712876595
$data2[14,1] = 73830.8632851734
$data2[14,2] = 73828.07639414954
$data2[14,3] = 73830.74848240045
$data2[14,4] = 73827.17590826476
$data2[14,5] = 73829.79154578631
$data2[14,6] = 73830.12926429899
$data2[14,7] = 73826.53576813792
$data2[14,8] = 73828.99813538382
$data2[14,9] = 73830.99813201187
$data2[15,0] = 6.090583347802331
$data2[15,1] = 6.096809103034817
$data2[15,2] = 6.094375516301245
$data2[15,3] = 6.09504423246972
$data2[15,4] = 6.092399715646361
$data2[15,5] = 6.092528228495246
$data2[15,6] = 6.094957554487915
$data2[15,7] = 6.095546975854501
$data2[15,8] = 6.098578235393386
$data2[15,9] = 6.097861709576327
$data2[16,0] = 87.66473788957597
$data2[16,1] = 88.10452110964768
$data2[16,2] = 87.69285093873057
$data2[16,3] = 88.13813779289703
$data2[16,4] = 88.29599320998808
$data2[16,5] = 88.37249115381557
$data2[16,6] = 87.79372876961028
$data2[16,7] = 87.9820387070034
$data2[16,8] = 88.25966293931403
$data2[16,9] = 88.18128768123452
$data2[17,0] = 50.00164019261909
$data2[17,1] = 49.99850087994426
$data2[17,2] = 49.99925892340833
$data2[17,3] = 49.9996273942193
$data2[17,4] = 50.0037324980015
$data2[17,5] = 50.00149993752392
$data2[17,6] = 49.9990622328394
$data2[17,7] = 49.99963059861378
$data2[17,8] = 49.99714851794409
$data2[17,9] = 50.00066299428491
$ws.Range("B16:K33").Value = $data2
